$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 184, shifting rows 184:245 down to 185:246.
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new data record.
$ws.Cells.Item(184, 1).Value = 8
$ws.Cells.Item(184, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 45141
$ws.Cells.Item(184, 5).Value = 4
$ws.Cells.Item(184, 6).Value = 100112044
$ws.Cells.Item(184, 7).Value = "Perejil"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 1200
$ws.Cells.Item(184, 11).Value = 2500
$ws.Cells.Item(184, 12).Value = 3000
$ws.Cells.Item(184, 13).Value = 2750
$ws.Cells.Item(184, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(184, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(184, 16).Value = 1833
$ws.Cells.Item(184, 17).Value = 1.5
$ws.Cells.Item(184, 18).Value = "Hortaliza"
